# Add 2022-Q4 data: new sheet inserted before the existing "2022-Q2" sheet,
# plus a new summary row on the "总计" sheet.
#
# NOTE: worksheet handles returned by Item()/Add() in this host resolve by
# tab *position*, so any handle obtained before a sheet is inserted/removed
# must be re-fetched (by name) before it is used again.

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) "总计" sheet: shift the existing 2022-Q2 summary row down to row 3 and
#    insert the new 2022-Q4 summary row in row 2.
# ---------------------------------------------------------------------------

$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("B3").NumberFormat = "@"
$totalSheet.Range("B3").Value = "2022-Q2"
$totalSheet.Range("B3").ClearFormats()
$totalSheet.Range("C3").Value = 23
$totalSheet.Range("D3").Value = 0.91

$totalSheet.Range("B2").NumberFormat = "@"
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("B2").ClearFormats()
$totalSheet.Range("C2").Value = 18
$totalSheet.Range("D2").Value = 1.07

# ---------------------------------------------------------------------------
# 2) New "2022-Q4" worksheet, inserted right before "2022-Q2", holding the
#    per-fund breakdown for the quarter.
# ---------------------------------------------------------------------------

$beforeSheet = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($beforeSheet)
$newSheet.Name = "2022-Q4"

# Re-fetch stable handles now that the sheet collection has changed.
$q4Sheet = $wb.Worksheets.Item("2022-Q4")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# Reuse the existing bold/bordered style (already used by the "总计" sheet)
# for the header row and the index column.
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2:A19").PasteSpecial(-4122)

$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$q4Rows = @(
  @("006803", "嘉实互通精选股票", "5.63", "89.43", "3.45", "0.1942", "10"),
  @("501099", "平安科技创新 3 年封闭混合", "2.73", "91.32", "4.71", "0.1286", "1"),
  @("011506", "建信高端装备股票A", "3.84", "89.26", "3.26", "0.1252", "6"),
  @("009008", "平安科技创新混合A", "2.49", "92.02", "4.92", "0.1225", "1"),
  @("011793", "建信智能汽车股票", "4.88", "87.61", "2.42", "0.1181", "7"),
  @("501098", "建信科技创新 3 年封闭混合", "3.25", "91.64", "3.51", "0.1141", "7"),
  @("001924", "华夏国企改革灵活配置混合", "2.49", "85.98", "3.42", "0.0852", "7"),
  @("009009", "平安科技创新混合C", "0.95", "92.02", "4.92", "0.0467", "1"),
  @("004265", "金鹰民丰回报定期开放混合", "4.49", "29.44", "0.67", "0.0301", "10"),
  @("011507", "建信高端装备股票C", "0.90", "89.26", "3.26", "0.0293", "6"),
  @("010571", "新沃创新领航混合C", "0.51", "93.56", "4.75", "0.0242", "4"),
  @("700004", "平安灵活配置混合A", "0.32", "79.07", "4.24", "0.0136", "2"),
  @("010570", "新沃创新领航混合A", "0.24", "93.56", "4.75", "0.0114", "4"),
  @("012143", "新沃内需增长混合A", "0.20", "93.63", "4.49", "0.0090", "4"),
  @("015078", "平安灵活配置混合C", "0.18", "79.07", "4.24", "0.0076", "2"),
  @("002564", "新沃通盈灵活配置混合", "0.10", "92.67", "6.30", "0.0063", "3"),
  @("501002", "长信价值优选混合", "0.40", "81.49", "1.25", "0.0050", "6"),
  @("012144", "新沃内需增长混合C", "0.04", "93.63", "4.49", "0.0018", "4")
)

$r = 2
foreach ($row in $q4Rows) {
  $q4Sheet.Cells.Item($r, 1).Value = $r - 2

  $q4Sheet.Cells.Item($r, 2).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 2).Value = $row[0]
  $q4Sheet.Cells.Item($r, 2).ClearFormats()

  $q4Sheet.Cells.Item($r, 3).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 3).Value = $row[1]
  $q4Sheet.Cells.Item($r, 3).ClearFormats()

  $q4Sheet.Cells.Item($r, 4).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 4).Value = $row[2]
  $q4Sheet.Cells.Item($r, 4).ClearFormats()

  $q4Sheet.Cells.Item($r, 5).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 5).Value = $row[3]
  $q4Sheet.Cells.Item($r, 5).ClearFormats()

  $q4Sheet.Cells.Item($r, 6).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 6).Value = $row[4]
  $q4Sheet.Cells.Item($r, 6).ClearFormats()

  $q4Sheet.Cells.Item($r, 7).NumberFormat = "@"
  $q4Sheet.Cells.Item($r, 7).Value = $row[5]
  $q4Sheet.Cells.Item($r, 7).ClearFormats()

  $q4Sheet.Cells.Item($r, 8).Value = [int]$row[6]

  $r = $r + 1
}

# Re-apply the bold index-column style that the ClearFormats() calls above
# stripped off A2:A19 (Value-only writes don't touch style, so this re-stamps
# the column once all the per-row text values have been written).
$totalSheet.Range("A2").Copy()
$q4Sheet.Range("A2:A19").PasteSpecial(-4122)
$r = 2
foreach ($row in $q4Rows) {
  $q4Sheet.Cells.Item($r, 1).Value = $r - 2
  $r = $r + 1
}

# ---------------------------------------------------------------------------
# 3) Restore "2022-Q2" as the selected/active sheet (unchanged from before).
# ---------------------------------------------------------------------------

$q2Sheet.Activate()
